# Reinstate v19 content and QA
# Resolves previously UNVERIFIED claims/figures with concrete source tags,
# per the "Cleanup pass" commit on 2026-02-08.

$wb = $excel.ActiveWorkbook

$claims  = $wb.Worksheets.Item("Claims")
$figures = $wb.Worksheets.Item("Figures")

# ---------------------------------------------------------------------------
# Claims sheet
# ---------------------------------------------------------------------------

# Row 67 (C066) - wallet capture / Pareto concentration guardrail
$claims.Range("C67").Value = 'Legacy v19 content on "wallet capture" and "Pareto concentration" is reintroduced as a strategic guardrail: a minority of high-spend households likely drives disproportionate preventive-care revenue concentration [S055; S054].'
$claims.Range("D67").Value = "S055, S054"
$claims.Range("E67").Value = "II.3 paragraph 2"
$claims.Range("J67").Value = "[2026-02-08] Cleanup pass: unresolved tail removed from section text; claim now fully source-tagged to S055/S054."

# Row 78 (C077) - Figure 38 TAM/SAM/SOM
$claims.Range("C78").Value = "Figure 38 TAM/SAM/SOM visual is mapped to in-repo figure data tab [S089, Tab: Figure 38]."
$claims.Range("D78").Value = "S089"
$claims.Range("J78").Value = "[2026-02-08] Cleanup pass: source linkage resolved to S089 Figure 38."

# Row 79 (C078) - FoodScience premium transaction
$claims.Range("C79").Value = "FoodScience appears in in-repo sponsor portfolio mapping (MSCP) and is used as transaction context in Part III table [S116, Tab: Sheet2]."
$claims.Range("D79").Value = "S116"
$claims.Range("E79").Value = "Table III.1 (FoodScience row)"
$claims.Range("J79").Value = "[2026-02-08] Cleanup pass: row rationale normalized to traceable portfolio mapping source."

# Row 80 (C079) - Figure 42 Margin Ladder
$claims.Range("C80").Value = "Figure 42 Margin Ladder visual is mapped to in-repo figure data tab [S089, Tab: Figure 44]."
$claims.Range("D80").Value = "S089"
$claims.Range("J80").Value = "[2026-02-08] Cleanup pass: source linkage resolved to S089 Figure 44."

# Row 81 (C080) - Figure 43 Strategic Capital Allocation matrix
$claims.Range("C81").Value = "Figure 43 Strategic Capital Allocation matrix is mapped to in-repo figure data tab [S089, Tab: Figure 45]."
$claims.Range("D81").Value = "S089"
$claims.Range("J81").Value = "[2026-02-08] Cleanup pass: source linkage resolved to S089 Figure 45."

# Row 84 (C083) - distribution-gatekeeper hypothesis
$claims.Range("C84").Value = "The v19 distribution-gatekeeper layer includes corporatized veterinary networks, specialty retail chains, and scaled e-commerce platforms (including IVC Evidensia, Zooplus, PetSmart, and Musti) [S116, Tab: Sheet1]."
$claims.Range("D84").Value = "S116"
$claims.Range("J84").Value = "[2026-02-08] Cleanup pass: converted prior unresolved hypothesis claim to sourced gatekeeper-mapping statement."

# ---------------------------------------------------------------------------
# Figures sheet
# ---------------------------------------------------------------------------

# Row 39 (FIG-38)
$figures.Range("D39").Value = "S089"
$figures.Range("I39").Value = "[2026-02-08] Cleanup pass: resolved from UNVERIFIED to S089 Figure 38."

# Row 43 (FIG-42)
$figures.Range("D43").Value = "S089"
$figures.Range("E43").Value = "Figure 44"
$figures.Range("I43").Value = "[2026-02-08] Cleanup pass: resolved from UNVERIFIED to S089 Figure 44 (margin ladder)."

# Row 44 (FIG-43)
$figures.Range("D44").Value = "S089"
$figures.Range("E44").Value = "Figure 45"
$figures.Range("I44").Value = "[2026-02-08] Cleanup pass: resolved from UNVERIFIED to S089 Figure 45 (strategic matrix)."
